$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Tnfsf11", "Tnfrsf11b", "ECs", 1, 0.3333333333333333, 0.05344666666666667, 0.16034, 0.01703322519986738, 0.01703322519986738, 1, 0.3333333333333333, 0.08241233333333334, 0.247237, 0.04727005612861496, 0.04727005612861496, 0.004404664508888889, 0.03964198058000001, 0.0008051615112490699, 0.00080516151124907),
    @("ECs", "Tnfsf11", "Tnfrsf11b", "FAPs", 1, 0.3333333333333333, 0.05344666666666667, 0.16034, 0.01703322519986738, 0.01703322519986738, 3, 1, 1.661024, 4.983072, 0.9527299438713851, 0.952729943871385, 0.08877619605333334, 0.79898576448, 0.01622806368861831, 0.01622806368861831),
    @("FAPs", "Tnfsf11", "Tnfrsf11b", "ECs", 3, 1, 3.017399, 9.052197, 0.9616322193748527, 0.9616322193748529, 1, 0.3333333333333333, 0.08241233333333334, 0.247237, 0.04727005612861496, 0.04727005612861496, 0.2486708921876667, 2.238038029689, 0.04545640898493386, 0.04545640898493387),
    @("FAPs", "Tnfsf11", "Tnfrsf11b", "FAPs", 3, 1, 3.017399, 9.052197, 0.9616322193748527, 0.9616322193748529, 3, 1, 1.661024, 4.983072, 0.9527299438713851, 0.952729943871385, 5.011972156575999, 45.107749409184, 0.9161758103899189, 0.9161758103899189),
    @("M1", "Tnfsf11", "Tnfrsf11b", "ECs", 1, 0.3333333333333333, 0.01182566666666667, 0.035477, 0.003768789637119215, 0.003768789637119216, 1, 0.3333333333333333, 0.08241233333333334, 0.247237, 0.04727005612861496, 0.04727005612861496, 0.0009745807832222223, 0.008771227049000001, 0.0001781508976835677, 0.0001781508976835678),
    @("M1", "Tnfsf11", "Tnfrsf11b", "FAPs", 1, 0.3333333333333333, 0.01182566666666667, 0.035477, 0.003768789637119215, 0.003768789637119216, 3, 1, 1.661024, 4.983072, 0.9527299438713851, 0.952729943871385, 0.01964271614933333, 0.176784445344, 0.003590638739435648, 0.003590638739435648),
    @("sCs", "Tnfsf11", "Tnfrsf11b", "ECs", 1, 0.3333333333333333, 0.05511766666666667, 0.165353, 0.0175657657881606, 0.0175657657881606, 1, 0.3333333333333333, 0.08241233333333334, 0.247237, 0.04727005612861496, 0.04727005612861496, 0.004542375517888889, 0.040881379661, 0.0008303347347484562, 0.0008303347347484562),
    @("sCs", "Tnfsf11", "Tnfrsf11b", "FAPs", 1, 0.3333333333333333, 0.05511766666666667, 0.165353, 0.0175657657881606, 0.0175657657881606, 3, 1, 1.661024, 4.983072, 0.9527299438713851, 0.952729943871385, 0.09155176715733335, 0.823965904416, 0.01673543105341215, 0.01673543105341214)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

Write-Host "Updated rows 2-9"